# Textbox response formatting fix
# Renames the five task-order sheets with refreshed timestamp suffixes and
# updates the stimulus-file names listed in column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (order: GNG, NB, RS, TOL, vSAT) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687280351787"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687314070442"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168731408046"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687314548554"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651168731516651"

# --- Sheet 1 (GNG) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "go_stims-16511687280042472.csv"
$ws.Range("B3").Value = "GNG_stims-16511687280196984.csv"
$ws.Range("B4").Value = "go_stims-16511687280206976.csv"
$ws.Range("B5").Value = "GNG_stims-16511687280351787.csv"

# --- Sheet 2 (NB) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = "TB-16511687313815978.csv"
$ws.Range("B3").Value = "OB-1651168729122648.csv"
$ws.Range("B4").Value = "OB-1651168729678065.csv"
$ws.Range("B5").Value = "TB-16511687303597918.csv"
$ws.Range("B6").Value = "ZB-match_3-1651168728784228.csv"
$ws.Range("B7").Value = "TB-16511687300293639.csv"
$ws.Range("B8").Value = "OB-1651168728997776.csv"
$ws.Range("B9").Value = "ZB-match_7-16511687282241187.csv"
$ws.Range("B10").Value = "ZB-match_1-16511687285271776.csv"

# --- Sheet 3 (RS) --- (name already updated above; no cell changes)

# --- Sheet 4 (TOL) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = "MM_stims-16511687314226515.csv"
$ws.Range("B3").Value = "ZM_stims-16511687314101226.csv"
$ws.Range("B4").Value = "MM_stims-1651168731438991.csv"
$ws.Range("B5").Value = "ZM_stims-16511687314226515.csv"
$ws.Range("B6").Value = "MM_stims-16511687314548554.csv"
$ws.Range("B7").Value = "ZM_stims-1651168731438991.csv"

# --- Sheet 5 (vSAT) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = "SAT_stims-1651168731470293.csv"
$ws.Range("B3").Value = "SAT_stims-16511687314578614.csv"
$ws.Range("B4").Value = "vSAT_stims-16511687314856663.csv"
$ws.Range("B5").Value = "vSAT_stims-16511687315014522.csv"
